# Lift state updating: progress, ButtonCallbackData to enum
#
# New people (Jukka, Juhani, Vladimir) are added to the "users" sheet and
# organised into groups ("groups" sheet) alongside the existing "CORE" crew
# (Eemeli, Akseli): a brand-new "LMG" group is introduced for them, and
# "Vallu" shows up as an extra tag in "users". Finally the "sites" sheet
# becomes the active/selected tab instead of "groups".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "users" sheet: the list of people moves from column A to column B, and
# three more names are appended; a "group" tag is added in column C.
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("users")

# Old data lived in A1:A2 ("Eemeli", "Akseli") - clear it out, it now
# lives one column over.
$wsUsers.Range("A1:A2").Clear()

$wsUsers.Range("B1").Value = "Eemeli"
$wsUsers.Range("B2").Value = "Akseli"
$wsUsers.Range("B3").Value = "Jukka"
$wsUsers.Range("B4").Value = "Juhani"
$wsUsers.Range("B5").Value = "Vladimir"
$wsUsers.Range("C5").Value = "Vallu"

$wsUsers.Range("B1:B5").Font.Bold = $true
$wsUsers.Range("B1:B5").HorizontalAlignment = -4108
$wsUsers.Range("B1:B5").VerticalAlignment = -4108

$wsUsers.Range("C5").HorizontalAlignment = -4108
$wsUsers.Range("C5").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# "groups" sheet: a second group "LMG" is added next to "CORE", with its
# own members in column B.
# ---------------------------------------------------------------------
$wsGroups = $wb.Worksheets.Item("groups")

$wsGroups.Range("B1").Value = "LMG"
$wsGroups.Range("B2").Value = "Jukka"
$wsGroups.Range("B3").Value = "Juhani"
$wsGroups.Range("B4").Value = "Vladimir"

$wsGroups.Range("B1").Font.Bold = $true
$wsGroups.Range("B1").HorizontalAlignment = -4108
$wsGroups.Range("B1").VerticalAlignment = -4108

$wsGroups.Range("B2:B4").HorizontalAlignment = -4108
$wsGroups.Range("B2:B4").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Window / selection state: "groups" and "users" keep a remembered
# selection but are no longer the active tab; "sites" becomes active.
# ---------------------------------------------------------------------
$wsGroups.Select()
$wsGroups.Range("E12").Select() | Out-Null

$wsUsers.Select()
$wsUsers.Range("F18").Select() | Out-Null

$wsSites = $wb.Worksheets.Item("sites")
$wsSites.Select()
$wsSites.Range("A18").Select() | Out-Null

# Best-effort: remembered application window bounds (maximized layout).
try {
    $excel.WindowState = -4137
    $excel.Left = -120
    $excel.Top = -120
    $excel.Width = 29040
    $excel.Height = 15990
} catch {
}
